$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and two rank-order swaps) to
# reflect the latest scrape, per the GitHub Actions "Updated cryptos list"
# commit. Numeric-looking price strings are written with a leading
# single-quote (quote-prefix) so Excel keeps them as text instead of
# silently coercing them to numbers and dropping significant trailing
# zeros (e.g. "1.00" -> 1, "3.70" -> 3.7).

$ws.Range("D2").Value = '42.090.55'
$ws.Range("E2").Value = '  +1.69%  '
$ws.Range("D3").Value = '2.216.00'
$ws.Range("E3").Value = '  +0.95%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''251.78'
$ws.Range("E5").Value = '  -1.80%  '
$ws.Range("D6").Value = '''0.623'
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("D7").Value = '''67.75'
$ws.Range("E7").Value = '  -1.01%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").Value = '''0.623'
$ws.Range("E9").Value = '  +6.86%  '
$ws.Range("D10").Value = '''38.93'
$ws.Range("E10").Value = '  +1.28%  '
$ws.Range("D11").Value = '''59.68'
$ws.Range("E11").Value = '  +2.67%  '
$ws.Range("D12").Value = '''0.0938'
$ws.Range("E12").Value = '  -0.79%  '
$ws.Range("D13").Value = '''7.04'
$ws.Range("E13").Value = '  -1.49%  '
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("D15").Value = '2.551.22'
$ws.Range("E15").Value = '  +1.25%  '
$ws.Range("D16").Value = '''0.878'
$ws.Range("E16").Value = '  +0.69%  '
$ws.Range("D17").Value = '''14.55'
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").Value = '2.219.81'
$ws.Range("E18").Value = '  +1.18%  '
$ws.Range("D19").Value = '42.003.20'
$ws.Range("E19").Value = '  +1.67%  '
$ws.Range("D20").Value = '0.0₃0960'
$ws.Range("E20").Value = '  +0.80%  '
$ws.Range("D21").Value = '''6.16'
$ws.Range("E21").Value = '  -1.31%  '
$ws.Range("D22").Value = '''72.47'
$ws.Range("E22").Value = '  +0.74%  '
$ws.Range("D23").Value = '''231.72'
$ws.Range("E23").Value = '  -0.59%  '
$ws.Range("D24").Value = '''2.02'
$ws.Range("E24").Value = '  -3.84%  '
$ws.Range("E25").Value = '  +1.06%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '''11.26'
$ws.Range("E27").Value = '  -6.96%  '
$ws.Range("E28").Value = '  -4.95%  '
$ws.Range("D29").Value = '''3.70'
$ws.Range("E29").Value = '  -1.47%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '''2.14'
$ws.Range("E30").Value = '  -2.09%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = '''166.70'
$ws.Range("E31").Value = '  -2.11%  '
$ws.Range("D32").Value = '''20.43'
$ws.Range("E32").Value = '  -1.19%  '
$ws.Range("D33").Value = '''5.93'
$ws.Range("E33").Value = '  +7.84%  '
$ws.Range("D34").Value = '''0.122'
$ws.Range("E34").Value = '  +1.10%  '
$ws.Range("D35").Value = '''0.0786'
$ws.Range("E35").Value = '  +7.78%  '
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").Value = '''4.60'
$ws.Range("E37").Value = '  -0.52%  '
$ws.Range("B38").Value = 'InjectiveProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D38").Value = '''25.87'
$ws.Range("E38").Value = '  +0.59%  '
$ws.Range("D39").Value = '''4.10'
$ws.Range("E39").Value = '  +2.95%  '
$ws.Range("E40").Value = '  +3.22%  '
$ws.Range("D41").Value = '''2.23'
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("D42").Value = '''12.03'
$ws.Range("E42").Value = '  -0.52%  '
$ws.Range("D43").Value = '''5.64'
$ws.Range("E43").Value = '  -2.39%  '
$ws.Range("D44").Value = '''5.07'
$ws.Range("E44").Value = '  +3.18%  '
$ws.Range("D45").Value = '''61.70'
$ws.Range("E45").Value = '  -4.63%  '
$ws.Range("D46").Value = '''0.196'
$ws.Range("E46").Value = '  -4.01%  '
$ws.Range("D47").Value = '''8.57'
$ws.Range("E47").Value = '  -0.66%  '
$ws.Range("D48").Value = '''0.0997'
$ws.Range("E48").Value = '  -2.23%  '
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("D50").Value = '''1.16'
$ws.Range("E50").Value = '  +0.75%  '
$ws.Range("D51").Value = '''4.34'
$ws.Range("E51").Value = '  +2.15%  '
